# "big chart atualizado com os testes" - update the Acceptance Tests ("Testes de
# Aceitação", column D) figure for the 06/10/2010 sample (row 4) from 7 to 13.
# This is the data point that feeds the big chart's third series.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("D4").Value = 13

# Leave the selection where the editor left it.
[void]$ws.Range("D5").Select()
